# The deck's slide-master theme ("Integral" palette, ppt/theme/theme1.xml) is
# recolored to the standard "Office Theme" palette (the palette that was
# previously only used by the Notes Master, ppt/theme/theme2.xml).
#
# PowerPoint's theme colors are edited through the 12-slot
# ThemeColorScheme on the master's Theme object (same order as the
# MsoThemeColorSchemeIndex enum: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink). Setting each slot's .RGB re-writes the <a:clrScheme> entries
# in the underlying theme part when the presentation is saved.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0         # dk1      -> 000000
$cs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477   # folHlink -> 954F72
